# The presentation ships two embedded theme parts:
#   ppt/theme/theme1.xml -> "Office Theme"  (currently only wired to the Notes Master)
#   ppt/theme/theme2.xml -> "Integral"      (wired to the Slide Master / the presentation's
#                                             live design, reachable as SlideMaster.Theme)
#
# The authored change swaps the content of those two theme parts, so the
# presentation's effective (visible) design becomes the plain "Office Theme"
# palette instead of "Integral".  Apply that by rewriting the live theme's
# 12-slot color scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) to the
# standard Office Theme RGB values.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$scheme = $master.Theme.ThemeColorScheme

$officeThemeRGB = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le $scheme.Count; $i++) {
    $scheme.Item($i).RGB = $officeThemeRGB[$i - 1]
}
